# Auto-generated edit script applying cell-value changes per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2576.25
$ws.Range("I62").Value = 2576.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2576.25
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1952.25
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2576.25
$ws.Range("I65").Value = 2576.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 12881.25
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -9761.25
$ws.Range("N65").ClearContents()

$ws.Range("H80").Value = 641.86365
$ws.Range("I80").Value = 586.5
$ws.Range("J80").Value = 708.3
$ws.Range("K80").Value = 1759.5
$ws.Range("L80").Value = 2124.9
$ws.Range("M80").Value = -761.5
$ws.Range("N80").Value = -4120.9

$ws.Range("H83").Value = 641.86365
$ws.Range("I83").Value = 586.5
$ws.Range("J83").Value = 708.3
$ws.Range("K83").Value = 5278.5
$ws.Range("L83").Value = 6374.7
$ws.Range("M83").Value = -286.5
$ws.Range("N83").Value = -16358.7

$ws.Range("H98").Value = 1056.3636
$ws.Range("I98").Value = 1110
$ws.Range("J98").Value = 716.6667
$ws.Range("K98").Value = 1110
$ws.Range("L98").Value = 716.6667
$ws.Range("M98").Value = 388
$ws.Range("N98").Value = -3712.6667

$ws.Range("H122").Value = 1056.3636
$ws.Range("I122").Value = 1110
$ws.Range("J122").Value = 716.6667
$ws.Range("K122").Value = 3330
$ws.Range("L122").Value = 2150.0001
$ws.Range("M122").Value = -880
$ws.Range("N122").Value = -7050.0001

$ws.Range("H125").Value = 6781.778
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 6781.778
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 61036.002
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -65956.00200000001

$ws.Range("H140").Value = 44000
$ws.Range("J140").Value = 44000
$ws.Range("L140").Value = 44000
$ws.Range("N140").Value = -54360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 726.36365
$ws.Range("I110").Value = 710.0714
$ws.Range("J110").Value = 754.875
$ws.Range("K110").Value = 710.0714
$ws.Range("L110").Value = 754.875
$ws.Range("M110").Value = 1334.9286
$ws.Range("N110").Value = -4844.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 12000
$ws.Range("J38").Value = 12000
$ws.Range("L38").Value = 12000
$ws.Range("N38").Value = -12832

$ws.Range("H99").Value = 125001736
$ws.Range("I99").Value = 166668060
$ws.Range("J99").Value = 2750
$ws.Range("K99").Value = 166668060
$ws.Range("L99").Value = 2750
$ws.Range("M99").Value = -166666562
$ws.Range("N99").Value = -5746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9014
$ws.Range("I31").Value = 803.4286
$ws.Range("J31").Value = 14761.4
$ws.Range("K31").Value = 803.4286
$ws.Range("L31").Value = 14761.4
$ws.Range("M31").Value = -508.4286
$ws.Range("N31").Value = -15351.4

$ws.Range("H34").Value = 9014
$ws.Range("I34").Value = 803.4286
$ws.Range("J34").Value = 14761.4
$ws.Range("K34").Value = 803.4286
$ws.Range("L34").Value = 14761.4
$ws.Range("M34").Value = -601.4286
$ws.Range("N34").Value = -15165.4

$ws.Range("H35").Value = 2404.2856
$ws.Range("I35").Value = 1600
$ws.Range("J35").Value = 7230
$ws.Range("K35").Value = 1600
$ws.Range("L35").Value = 7230
$ws.Range("M35").Value = -1306
$ws.Range("N35").Value = -7818

$ws.Range("H99").Value = 3559.3333
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 5912.857
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 5912.857
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -8908.857

$ws.Range("H126").Value = 3559.3333
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 5912.857
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 17738.571
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -22678.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3613.5293
$ws.Range("I125").Value = 1466.6666
$ws.Range("J125").Value = 4073.5715
$ws.Range("K125").Value = 4399.9998
$ws.Range("L125").Value = 12220.7145
$ws.Range("M125").Value = 520.0002000000004
$ws.Range("N125").Value = -22060.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1692.2778
$ws.Range("I102").Value = 1747.4
$ws.Range("J102").Value = 1416.6666
$ws.Range("K102").Value = 1747.4
$ws.Range("L102").Value = 1416.6666
$ws.Range("M102").Value = -125.4000000000001
$ws.Range("N102").Value = -4660.6666

$ws.Range("H122").Value = 2947938.8
$ws.Range("I122").Value = 4322671
$ws.Range("J122").Value = 2084.1428
$ws.Range("K122").Value = 12968013
$ws.Range("L122").Value = 6252.428400000001
$ws.Range("M122").Value = -12965563
$ws.Range("N122").Value = -11152.4284

$ws.Range("H126").Value = 8613.933999999999
$ws.Range("I126").Value = 10828.546
$ws.Range("J126").Value = 2523.75
$ws.Range("K126").Value = 32485.638
$ws.Range("L126").Value = 7571.25
$ws.Range("M126").Value = -30015.638
$ws.Range("N126").Value = -12511.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 90911530
$ws.Range("I40").Value = 100001784
$ws.Range("J40").Value = 9005
$ws.Range("K40").Value = 100001784
$ws.Range("L40").Value = 9005
$ws.Range("M40").Value = -100001648
$ws.Range("N40").Value = -9277

$ws.Range("H122").Value = 7939575
$ws.Range("I122").Value = 8931647
$ws.Range("K122").Value = 26794941
$ws.Range("M122").Value = -26792491

$ws.Range("H127").Value = 89666
$ws.Range("J127").Value = 89666
$ws.Range("L127").Value = 89666
$ws.Range("N127").Value = -99586

$ws.Range("H132").Value = 12386319
$ws.Range("I132").Value = 16672825
$ws.Range("J132").Value = 3078.6667
$ws.Range("K132").Value = 50018475
$ws.Range("L132").Value = 9236.000100000001
$ws.Range("M132").Value = -50015945
$ws.Range("N132").Value = -14296.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 40000
$ws.Range("I109").Value = 40000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 40000
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -38613
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 1855
$ws.Range("I122").Value = 2066.6667
$ws.Range("J122").Value = 1537.5
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("L122").Value = 4612.5
$ws.Range("M122").Value = -3750.000100000001
$ws.Range("N122").Value = -9512.5

